$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds text-valued period labels (e.g. "01-04-2021"). Excel's COM
# layer auto-parses strings that look like dates into date serials, so force
# the cell to Text format before entering the value, then restore the
# "Normal" style afterwards so the cell keeps the workbook's default
# (unstyled) look, matching the rest of the column.
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "01-07-2021"
$ws.Range("A76").Style = "Normal"
$ws.Range("B76").Value = 21658
$ws.Range("C76").Value = 1659
$ws.Range("D76").Value = -748
$ws.Range("E76").Value = 1619
$ws.Range("F76").Value = 788
$ws.Range("G76").Value = 16348
$ws.Range("H76").Value = 663
$ws.Range("I76").Value = 15685
$ws.Range("J76").Value = 15437
$ws.Range("K76").Value = 248
$ws.Range("L76").Value = -2635
$ws.Range("M76").Value = 6286
$ws.Range("N76").Value = 1116
$ws.Range("O76").Value = 2712
$ws.Range("P76").Value = 80
$ws.Range("Q76").Value = 0
$ws.Range("R76").Value = 2378
